$d = $word.ActiveDocument

# 1) Table edits: UI test table (first table), rows 6-8 (1-indexed), columns 2 and 3.
$tbl = $d.Tables.Item(1)

$tbl.Cell(6, 2).Range.Text = "Disable the UI"
$tbl.Cell(6, 3).Range.Text = "Appropriate UI elements are disables"

$tbl.Cell(7, 2).Range.Text = "Enable UI"
$tbl.Cell(7, 3).Range.Text = "Appropriate elements are enabled"

$tbl.Cell(8, 2).Range.Text = ""
$tbl.Cell(8, 3).Range.Text = ""

# 2) Move the _GoBack bookmark: remove old, add new at end of cell (7,3) text.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$cell = $tbl.Cell(7, 3)
$cellRange = $cell.Range
$bmRange = $d.Range($cellRange.End - 2, $cellRange.End - 2)
$d.Bookmarks.Add("_GoBack", $bmRange)
